$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for Case_3_58 (380 kV case)
$rowCount = 24

$BData = New-Object 'object[,]' $rowCount,1
$BData[0,0] = 7.197104432957967
$BData[1,0] = 6.985711345888937
$BData[2,0] = 6.854309205985398
$BData[3,0] = 6.800451048557073
$BData[4,0] = 6.791492008661015
$BData[5,0] = 6.853583986876178
$BData[6,0] = 7.124603458397402
$BData[7,0] = 7.639543133442873
$BData[8,0] = 8.003171866988232
$BData[9,0] = 8.164558602311233
$BData[10,0] = 8.225030390514073
$BData[11,0] = 8.212036183194444
$BData[12,0] = 8.169546851824519
$BData[13,0] = 8.143435574808754
$BData[14,0] = 7.992538298014246
$BData[15,0] = 7.898889186358264
$BData[16,0] = 7.844648711769666
$BData[17,0] = 7.826221247817052
$BData[18,0] = 7.908897678758637
$BData[19,0] = 8.182044884703297
$BData[20,0] = 8.35679649208698
$BData[21,0] = 8.263891860112842
$BData[22,0] = 7.904374081553639
$BData[23,0] = 7.502510398995542
$ws.Range("B2:B25").Value = $BData

$CData = New-Object 'object[,]' $rowCount,1
$CData[0,0] = 4.46485883537922
$CData[1,0] = 4.292404538252832
$CData[2,0] = 4.181756748038818
$CData[3,0] = 4.135503674997983
$CData[4,0] = 4.127754210531358
$CData[5,0] = 4.181137624699352
$CData[6,0] = 4.406403437414089
$CData[7,0] = 4.809117406121309
$CData[8,0] = 5.079784655031846
$CData[9,0] = 5.197215263212902
$CData[10,0] = 5.240848490193974
$CData[11,0] = 5.231488666703109
$CData[12,0] = 5.20082181595845
$CData[13,0] = 5.181928287452101
$CData[14,0] = 5.071994224189761
$CData[15,0] = 5.003081866167844
$CData[16,0] = 4.96290982364664
$CData[17,0] = 4.949216819394458
$CData[18,0] = 5.010473247637387
$CData[19,0] = 5.209852194891969
$CData[20,0] = 5.335284024826882
$CData[21,0] = 5.268789280658707
$CData[22,0] = 5.007133326588125
$CData[23,0] = 4.704510895555771
$ws.Range("C2:C25").Value = $CData

$DData = New-Object 'object[,]' $rowCount,1
$DData[0,0] = 5.176866307448251
$DData[1,0] = 5.155094932509346
$DData[2,0] = 5.141252052318277
$DData[3,0] = 5.135492694787291
$DData[4,0] = 5.134529250376517
$DData[5,0] = 5.141174856872444
$DData[6,0] = 5.169458877885956
$DData[7,0] = 5.221102107091446
$DData[8,0] = 5.256648362555091
$DData[9,0] = 5.272285356577234
$DData[10,0] = 5.278128849427157
$DData[11,0] = 5.276873832517097
$DData[12,0] = 5.272767672392406
$DData[13,0] = 5.270242350119385
$DData[14,0] = 5.255615613980821
$DData[15,0] = 5.246505146357173
$DData[16,0] = 5.241214933430154
$DData[17,0] = 5.239415189007198
$DData[18,0] = 5.247480166638711
$DData[19,0] = 5.273975874949437
$DData[20,0] = 5.290837750936875
$DData[21,0] = 5.281880246720212
$DData[22,0] = 5.247039523351053
$DData[23,0] = 5.207547291825257
$ws.Range("D2:D25").Value = $DData

$FData = New-Object 'object[,]' $rowCount,1
$FData[0,0] = 27.10163829750831
$FData[1,0] = 27.0149316256284
$FData[2,0] = 26.96837216141881
$FData[3,0] = 26.9510885198841
$FData[4,0] = 26.94832089943381
$FData[5,0] = 26.96813221446597
$FData[6,0] = 27.07036243387794
$FData[7,0] = 27.32328605414124
$FData[8,0] = 27.54022899827629
$FData[9,0] = 27.64545711226965
$FData[10,0] = 27.68622300551603
$FData[11,0] = 27.6774028696845
$FData[12,0] = 27.64879267556786
$FData[13,0] = 27.63138700904812
$FData[14,0] = 27.53348183781063
$FData[15,0] = 27.47507997237989
$FData[16,0] = 27.44210574542144
$FData[17,0] = 27.4310478513082
$FData[18,0] = 27.48123325017512
$FData[19,0] = 27.65717144550898
$FData[20,0] = 27.77749830553574
$FData[21,0] = 27.71279662014011
$FData[22,0] = 27.47844947720704
$FData[23,0] = 27.24932883013334
$ws.Range("F2:F25").Value = $FData

$GData = New-Object 'object[,]' $rowCount,1
$GData[0,0] = 33.87476623968679
$GData[1,0] = 33.69133336912326
$GData[2,0] = 33.58801092124315
$GData[3,0] = 33.54827464489092
$GData[4,0] = 33.54182026475845
$GData[5,0] = 33.58746539958569
$GData[6,0] = 33.80960535134113
$GData[7,0] = 34.31762699956496
$GData[8,0] = 34.7328364922998
$GData[9,0] = 34.93030350819038
$GData[10,0] = 35.00626413875828
$GData[11,0] = 34.98985282397516
$GData[12,0] = 34.93652940130963
$GData[13,0] = 34.90401993355135
$GData[14,0] = 34.72010007570594
$GData[15,0] = 34.60943522909093
$GData[16,0] = 34.54659510866438
$GData[17,0] = 34.5254593728043
$GData[18,0] = 34.62113207453174
$GData[19,0] = 34.95216006090735
$GData[20,0] = 35.17538296966131
$GData[21,0] = 35.05563280122284
$GData[22,0] = 34.61584148936517
$GData[23,0] = 34.17265247234582
$ws.Range("G2:G25").Value = $GData

$HData = New-Object 'object[,]' $rowCount,1
$HData[0,0] = 15.34958577672841
$HData[1,0] = 15.36885259539231
$HData[2,0] = 15.38323650227154
$HData[3,0] = 15.38973922278377
$HData[4,0] = 15.39085768552462
$HData[5,0] = 15.38332160567415
$HData[6,0] = 15.35569828860355
$HData[7,0] = 15.3218376207024
$HData[8,0] = 15.30939493708051
$HData[9,0] = 15.30644331715996
$HData[10,0] = 15.30571553319001
$HData[11,0] = 15.30585492699885
$HData[12,0] = 15.30637562523216
$HData[13,0] = 15.30674535854185
$HData[14,0] = 15.30964237748603
$HData[15,0] = 15.31211370653459
$HData[16,0] = 15.31379008184377
$HData[17,0] = 15.31440144268831
$HData[18,0] = 15.31182424117999
$HData[19,0] = 15.30621209826646
$HData[20,0] = 15.30481713952665
$HData[21,0] = 15.3053535825244
$HData[22,0] = 15.31195431244137
$HData[23,0] = 15.32881749702926
$ws.Range("H2:H25").Value = $HData

$IData = New-Object 'object[,]' $rowCount,1
$IData[0,0] = 22.59092971433252
$IData[1,0] = 22.61230207878932
$IData[2,0] = 22.62973177519134
$IData[3,0] = 22.63791477924383
$IData[4,0] = 22.63933871867124
$IData[5,0] = 22.62983776396961
$IData[6,0] = 22.5974034228609
$IData[7,0] = 22.56808564390244
$IData[8,0] = 22.5675832831748
$IData[9,0] = 22.57194459722898
$IData[10,0] = 22.57425718463186
$IData[11,0] = 22.57372971476891
$IData[12,0] = 22.57212160219105
$IData[13,0] = 22.57122269932458
$IData[14,0] = 22.56739071830524
$IData[15,0] = 22.56621635613939
$IData[16,0] = 22.56597287481537
$IData[17,0] = 22.56596459105658
$IData[18,0] = 22.56629665184927
$IData[19,0] = 22.57257599676448
$IData[20,0] = 22.58053322903153
$IData[21,0] = 22.57593348605894
$IData[22,0] = 22.56625900567061
$IData[23,0] = 22.57233086675245
$ws.Range("I2:I25").Value = $IData

$KData = New-Object 'object[,]' $rowCount,1
$KData[0,0] = 7.493077718613511
$KData[1,0] = 7.363290461200854
$KData[2,0] = 7.284192867791535
$KData[3,0] = 7.252155661876677
$KData[4,0] = 7.24684910083228
$KData[5,0] = 7.283759948397756
$KData[6,0] = 7.448232073658547
$KData[7,0] = 7.773433052663369
$KData[8,0] = 8.056749574135678
$KData[9,0] = 8.222713787200737
$KData[10,0] = 8.284446950136532
$KData[11,0] = 8.2712015295547
$KData[12,0] = 8.227815073281915
$KData[13,0] = 8.201093837444306
$KData[14,0] = 8.045748897359594
$KData[15,0] = 7.948492747760034
$KData[16,0] = 7.90686867269154
$KData[17,0] = 7.894798605650958
$KData[18,0] = 7.958919577382479
$KData[19,0] = 8.240589149249434
$KData[20,0] = 8.418174177228286
$KData[21,0] = 8.32399644085141
$KData[22,0] = 7.954207897902032
$KData[23,0] = 7.685450123647988
$ws.Range("K2:K25").Value = $KData

$MData = New-Object 'object[,]' $rowCount,1
$MData[0,0] = 19.97696696572945
$MData[1,0] = 19.38610794789741
$MData[2,0] = 19.02231616688643
$MData[3,0] = 18.87406365757266
$MData[4,0] = 18.8494536277925
$MData[5,0] = 19.02031646411048
$MData[6,0] = 19.77359157351269
$MData[7,0] = 21.23258364199129
$MData[8,0] = 22.28077389036368
$MData[9,0] = 22.75011847002204
$MData[10,0] = 22.92659797367991
$MData[11,0] = 22.88864820711136
$MData[12,0] = 22.76466356202744
$MData[13,0] = 22.68855154738805
$MData[14,0] = 22.24993519435257
$MData[15,0] = 21.97881183934767
$MData[16,0] = 21.82217396283972
$MData[17,0] = 21.76902533517152
$MData[18,0] = 22.00774662469286
$MData[19,0] = 22.8011161242409
$MData[20,0] = 23.31226718440264
$MData[21,0] = 23.04018335597244
$MData[22,0] = 21.99466757896405
$MData[23,0] = 20.84120100499489
$ws.Range("M2:M25").Value = $MData

